$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer (primary, odd pages) - Pearson Edexcel logo: image2.png -> image1.png
$ftr1 = $sec.Footers.Item(1)
if ($ftr1.Range.InlineShapes.Count -ge 1) {
    $ftr1.Range.InlineShapes.Item(1).Name = "image1.png"
}

# Footer (even pages) - Pearson Edexcel logo: image2.png -> image1.png
$ftr2 = $sec.Footers.Item(2)
if ($ftr2.Range.InlineShapes.Count -ge 1) {
    $ftr2.Range.InlineShapes.Item(1).Name = "image1.png"
}

# Header (even pages) - BTEC logo: image1.jpg -> image2.jpg
$hdr2 = $sec.Headers.Item(2)
if ($hdr2.Range.InlineShapes.Count -ge 1) {
    $hdr2.Range.InlineShapes.Item(1).Name = "image2.jpg"
}

Write-Output "Renamed inline picture shapes."
